$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Type Group" header column, matching the style of the
# existing header row (row 2), and select it as the active cell.
$ws.Range("K2").Copy()
$ws.Range("L2").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("L2").Value = "Type Group"
$ws.Range("L2").Select()
